$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prime a bold font into the styles table (without leaving any cell bold) ---
# Touching Font.Bold on a cell causes the workbook to register a bold "Calibri 11"
# font entry in xl/styles.xml; flipping it back off keeps the cell itself on the
# default (non-bold) style while the bold font stays registered for reuse by the
# rich-text runs we add below.
$primeCell = $ws.Cells.Item(28, 2)
$primeCell.Font.Bold = $true
$primeCell.Font.Bold = $false

# --- Row 33 (git diff): clarify the description text ---
$ws.Cells.Item(33, 2).Value = "Shows difference between working directory copy and git repository copy"

# --- Row 28 (git push origin master): replace plain text with rich text ---
# that calls out "origin" and "master" in bold.
$cell = $ws.Cells.Item(28, 2)
$text = "update remote repository with local (origin -> Name of remote repo; master -> branch on remote repo)"
$cell.Value = $text
$cell.Characters(38, 6).Font.Bold = $true
$cell.Characters(69, 6).Font.Bold = $true

# --- Widen column B to fit the new, longer text ---
$ws.Columns.Item(2).ColumnWidth = 93.65

# --- Update the selection / active cell ---
$ws.Range("B30").Select() | Out-Null

# --- Switch the page to portrait orientation ---
$ws.PageSetup.Orientation = 1
